$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" "256.79"
Set-TextValue $ws "E2" "-0.75%"
Set-TextValue $ws "G2" "16"

Set-TextValue $ws "D3" "27.09"
Set-TextValue $ws "E3" "-3.44%"
Set-TextValue $ws "G3" "16"

Set-TextValue $ws "D4" "4.626"
Set-TextValue $ws "E4" "-11.30%"
Set-TextValue $ws "G4" "16"

Set-TextValue $ws "D5" "0.05894"
Set-TextValue $ws "E5" "-0.19%"
Set-TextValue $ws "G5" "16"

Set-TextValue $ws "D6" "6.646"
Set-TextValue $ws "E6" "-0.87%"
Set-TextValue $ws "G6" "16"

Set-TextValue $ws "D7" "0.8646"
Set-TextValue $ws "E7" "-0.36%"
Set-TextValue $ws "G7" "16"

Set-TextValue $ws "D8" "0.9333"
Set-TextValue $ws "E8" "-8.63%"
Set-TextValue $ws "G8" "16"

Set-TextValue $ws "D9" "0.1403"
Set-TextValue $ws "E9" "-0.81%"
Set-TextValue $ws "G9" "16"

Set-TextValue $ws "D10" "0.03815"
Set-TextValue $ws "E10" "8.27%"
Set-TextValue $ws "G10" "16"

Set-TextValue $ws "E11" "-1.00%"
Set-TextValue $ws "G11" "16"

Set-TextValue $ws "D12" "0.03197"
Set-TextValue $ws "E12" "1.44%"
Set-TextValue $ws "G12" "16"

Set-TextValue $ws "D13" "0.09225"
Set-TextValue $ws "E13" "0.11%"
Set-TextValue $ws "G13" "16"

Set-TextValue $ws "D14" "0.001556"
Set-TextValue $ws "E14" "0.97%"
Set-TextValue $ws "G14" "16"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D15" "0.006106"
Set-TextValue $ws "E15" "3.74%"
Set-TextValue $ws "G15" "16"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D16" "3.515"
Set-TextValue $ws "E16" "0.39%"
Set-TextValue $ws "G16" "16"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D17" "3.192"
Set-TextValue $ws "E17" "-1.07%"
Set-TextValue $ws "G17" "16"

$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D18" "2.212"
Set-TextValue $ws "E18" "0.34%"
Set-TextValue $ws "G18" "16"

$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D19" "0.01039"
Set-TextValue $ws "E19" "1,610.22%"
Set-TextValue $ws "G19" "16"

Set-TextValue $ws "D20" "0.3118"
Set-TextValue $ws "E20" "-1.82%"
Set-TextValue $ws "G20" "16"

Set-TextValue $ws "E21" "-1.10%"
Set-TextValue $ws "G21" "16"

Set-TextValue $ws "D22" "3.846"
Set-TextValue $ws "E22" "8.22%"
Set-TextValue $ws "G22" "16"

Set-TextValue $ws "D23" "0.04210"
Set-TextValue $ws "E23" "0.31%"
Set-TextValue $ws "G23" "16"

Set-TextValue $ws "D24" "0.001219"
Set-TextValue $ws "E24" "-0.19%"
Set-TextValue $ws "G24" "16"

Set-TextValue $ws "D25" "0.004282"
Set-TextValue $ws "E25" "-6.14%"
Set-TextValue $ws "G25" "16"

Set-TextValue $ws "E26" "0.16%"
Set-TextValue $ws "G26" "16"

Set-TextValue $ws "D27" "0.0001936"
Set-TextValue $ws "E27" "31.75%"
Set-TextValue $ws "G27" "16"

Set-TextValue $ws "G28" "16"

Set-TextValue $ws "G29" "16"

Set-TextValue $ws "G30" "16"

Set-TextValue $ws "G31" "16"

Set-TextValue $ws "G32" "16"

Set-TextValue $ws "G33" "16"

Set-TextValue $ws "G34" "16"

Set-TextValue $ws "G35" "16"

Set-TextValue $ws "G36" "16"

Set-TextValue $ws "G37" "16"

Set-TextValue $ws "G38" "16"

Set-TextValue $ws "G39" "16"

Set-TextValue $ws "D40" "0.03827"
Set-TextValue $ws "E40" "-0.30%"
Set-TextValue $ws "G40" "16"

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006166"
Set-TextValue $ws "E41" "12.98%"
Set-TextValue $ws "G41" "16"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1098"
Set-TextValue $ws "E42" "-0.35%"
Set-TextValue $ws "G42" "16"

Set-TextValue $ws "D43" "0.002263"
Set-TextValue $ws "E43" "-4.52%"
Set-TextValue $ws "G43" "16"

Set-TextValue $ws "D44" "0.01133"
Set-TextValue $ws "E44" "19.11%"
Set-TextValue $ws "G44" "16"

Set-TextValue $ws "D45" "0.00005474"
Set-TextValue $ws "E45" "0.97%"
Set-TextValue $ws "G45" "16"

Set-TextValue $ws "D46" "0.00000000750"
Set-TextValue $ws "E46" "0.15%"
Set-TextValue $ws "G46" "16"

Set-TextValue $ws "D47" "0.06017"
Set-TextValue $ws "E47" "-36.53%"
Set-TextValue $ws "G47" "16"

Set-TextValue $ws "G48" "16"

Set-TextValue $ws "D49" "0.00002099"
Set-TextValue $ws "E49" "0.15%"
Set-TextValue $ws "G49" "16"

Set-TextValue $ws "D50" "0.0001999"
Set-TextValue $ws "E50" "0.15%"
Set-TextValue $ws "G50" "16"

Set-TextValue $ws "G51" "16"
